# Depreciation calculator and analysis files
#
# The "Information" sheet (a cover sheet holding only a text-box drawing
# with project notes) is removed from the workbook, leaving
# "Depreciation Calculator" as the sole worksheet.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$infoSheet = $wb.Worksheets.Item("Information")
if ($infoSheet -ne $null) {
    $infoSheet.Delete()
}

# Make the remaining sheet the active one so the saved workbook view
# doesn't keep pointing at a now-nonexistent tab index.
$wb.Worksheets.Item("Depreciation Calculator").Activate()

$excel.DisplayAlerts = $true
